$wb = $excel.ActiveWorkbook

$wsSweetener = $wb.Worksheets.Item("sweetener")
$wsSugar = $wb.Worksheets.Item("sugar")
$wsArtificial = $wb.Worksheets.Item("artificial sweetener")

# Fix the "food supergroup" label -> "food super group" on every sheet
# (sweetener: row 7, sugar/artificial sweetener: row 3)
$wsSweetener.Range("A7").Value = "food super group"
$wsSugar.Range("A3").Value = "food super group"
$wsArtificial.Range("A3").Value = "food super group"

# Update each sheet's selection/active cell to rest on that label cell
$wsSugar.Activate()
$wsSugar.Range("A3").Select()

$wsArtificial.Activate()
$wsArtificial.Range("A3").Select()

# Finally, make "sweetener" the active tab/sheet with A7 selected,
# matching the workbook's saved view state
$wsSweetener.Activate()
$wsSweetener.Range("A7").Select()
